# Refresh the crypto price/volume table on Sheet1 (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.418.28'
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").Value = '2.331.20'
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("E4").Value = '  -0.43%  '
$c = $ws.Range("D5")
$c.Value = '''512.80'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.86%  '
$c = $ws.Range("D6")
$c.Value = '''132.75'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  -0.36%  '
$c = $ws.Range("D9")
$c.Value = '''0.100'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.51%  '
$ws.Range("E10").Value = '  -1.01%  '
$c = $ws.Range("D11")
$c.Value = '''5.29'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.69%  '
$c = $ws.Range("D12")
$c.Value = '''0.338'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.69%  '
$c = $ws.Range("D13")
$c.Value = '''23.61'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.42%  '
$ws.Range("D14").Value = '2.748.36'
$ws.Range("E14").Value = '  -0.38%  '
$ws.Range("D15").Value = '56.380.33'
$ws.Range("E15").Value = '  -0.53%  '
$c = $ws.Range("D16")
$c.Value = '''0.0000132'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.86%  '
$ws.Range("D17").Value = '2.345.13'
$ws.Range("E17").Value = '  +0.20%  '
$c = $ws.Range("D18")
$c.Value = '''10.39'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.04%  '
$c = $ws.Range("D19")
$c.Value = '''323.71'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.95%  '
$c = $ws.Range("D20")
$c.Value = '''4.16'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.96%  '
$c = $ws.Range("D21")
$c.Value = '''6.63'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.91%  '
$c = $ws.Range("D22")
$c.Value = '''1.00'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.01%  '
$c = $ws.Range("D23")
$c.Value = '''61.39'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.07%  '
$c = $ws.Range("D24")
$c.Value = '''8.61'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +11.65%  '
$ws.Range("E25").Value = '  +2.44%  '
$c = $ws.Range("D26")
$c.Value = '''1.00'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("E27").Value = '  +6.70%  '
$c = $ws.Range("D28")
$c.Value = '''167.45'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -2.02%  '
$ws.Range("D29").Value = '0.0₃0720'
$ws.Range("E29").Value = '  -1.99%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("E31").Value = '  -1.47%  '
$c = $ws.Range("D32")
$c.Value = '''18.31'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.38%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("E34").Value = '  +0.21%  '
$c = $ws.Range("D35")
$c.Value = '''1.26'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +2.12%  '
$c = $ws.Range("D36")
$c.Value = '''3.95'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("E37").Value = '  -5.82%  '
$ws.Range("E38").Value = '  +2.33%  '
$c = $ws.Range("D39")
$c.Value = '''38.43'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.53%  '
$c = $ws.Range("D40")
$c.Value = '''150.35'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +8.98%  '
$c = $ws.Range("D41")
$c.Value = '''0.374'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.16%  '
$c = $ws.Range("D42")
$c.Value = '''3.56'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.16%  '
$c = $ws.Range("D43")
$c.Value = '''277.88'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.86%  '
$c = $ws.Range("D44")
$c.Value = '''5.08'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +1.11%  '
$c = $ws.Range("D45")
$c.Value = '''0.0925'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.47%  '
$ws.Range("E46").Value = '  -1.04%  '
$c = $ws.Range("D47")
$c.Value = '''0.554'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.68%  '
$c = $ws.Range("D48")
$c.Value = '''18.13'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +6.23%  '
$ws.Range("E49").Value = '  +0.29%  '
$c = $ws.Range("D50")
$c.Value = '''0.0214'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.63%  '
$c = $ws.Range("D51")
$c.Value = '''17.02'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +2.01%  '
